# Insert a new data row at row 322 (pushing existing rows 322-397 down to 323-398)
# and populate it with the new weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 322; this shifts rows 322:397 down to 323:398
$ws.Rows.Item(322).Insert()

# Populate the newly inserted row 322 with the new record
$ws.Range("A322").Value = 6
$ws.Range("B322").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C322").Value = "Metropolitana"
$ws.Range("D322").Value = 44543
$ws.Range("D322").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E322").Value = 13
$ws.Range("F322").Value = 100112044
$ws.Range("G322").Value = "Perejil"
$ws.Range("H322").Value = "Sin especificar"
$ws.Range("I322").Value = "Primera"
$ws.Range("J322").Value = 150
$ws.Range("K322").Value = 9500
$ws.Range("L322").Value = 10000
$ws.Range("M322").Value = 9700
$ws.Range("N322").Value = '$/docena de atados'
$ws.Range("O322").Value = "Región Metropolitana"
$ws.Range("P322").Value = 3233
$ws.Range("Q322").Value = 3
$ws.Range("R322").Value = "Hortaliza"
